$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1096
$ws1.Range("F3").Value = 4740
$ws1.Range("F4").Value = 621
$ws1.Range("F6").Value = 1904
$ws1.Range("F15").Value = 4
$ws1.Range("F16").Value = 1931
$ws1.Range("F17").Value = 597
$ws1.Range("F18").Value = 9
$ws1.Range("F19").Value = 535
$ws1.Range("F21").Value = 223
$ws1.Range("F22").Value = 75
$ws1.Range("F23").Value = 75
$ws1.Range("F27").Value = 2537
$ws1.Range("F28").Value = 19
$ws1.Range("F29").Value = 10
$ws1.Range("F31").Value = 1634
$ws1.Range("F36").Value = 4342

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 294
$ws2.Range("F17").Value = 293
$ws2.Range("F28").Value = 7
$ws2.Range("F29").Value = 94

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1749
$ws3.Range("F6").Value = 1096
$ws3.Range("F7").Value = 380

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1749
$ws4.Range("F4").Value = 1096
$ws4.Range("F5").Value = 380
$ws4.Range("F7").Value = 1096
$ws4.Range("F9").Value = 4740
$ws4.Range("F10").Value = 621
$ws4.Range("F12").Value = 1904
$ws4.Range("F24").Value = 4
$ws4.Range("F25").Value = 1931
$ws4.Range("F26").Value = 597
$ws4.Range("F27").Value = 9
$ws4.Range("F28").Value = 535
$ws4.Range("F29").Value = 223
$ws4.Range("F30").Value = 75
$ws4.Range("F31").Value = 294
$ws4.Range("F32").Value = 293
$ws4.Range("F39").Value = 2537
$ws4.Range("F45").Value = 1634
$ws4.Range("F49").Value = 4342
